# Gantt chart update: move the Display Week forward and stretch/re-order a
# few task durations on the ProjectSchedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
$ws.Activate()

# Display Week: 1 -> 10 (drives the whole I4:CG5 date header via formulas)
$ws.Range("E4").Value = 10

# Task 1 (row 9): end date pushed out, was +6 days, now +14 days
$ws.Range("F9").Formula = "=E9+14"

# Task 2 (row 10): end date pushed out, was +6 days, now +20 days
$ws.Range("F10").Formula = "=E10+20"

# Phase 2 sub-tasks (rows 12 & 13): swap the two task names...
$b12 = $ws.Range("B12").Value()
$b13 = $ws.Range("B13").Value()
$ws.Range("B12").Value = $b13
$ws.Range("B13").Value = $b12

# ...and re-point/extend their start & end date formulas
$ws.Range("E12").Formula = "=F9+1"
$ws.Range("F12").Formula = "=E12+20"
$ws.Range("F13").Formula = "=E13+6"

# Rows 16, 18, 21 & 23: extend end dates from +6 days to +13 days
$ws.Range("F16").Formula = "=E16+13"
$ws.Range("F18").Formula = "=E18+13"
$ws.Range("F21").Formula = "=E21+13"
$ws.Range("F23").Formula = "=E23+13"

# Leave the cursor roughly where the author left it
$ws.Range("F4").Select()
